$wb = $excel.ActiveWorkbook

$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")

$hotelInfo.Range("C1").EntireColumn.Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"

$reviewInfo.Move($hotelInfo)
